$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RealDevices")
$ws2 = $wb.Worksheets.Item("Browsers")

# --- RealDevices (sheet1): new "supported" column (L) + platformVersion fills ---
$ws1.Range("L1").Value = "supported"
$ws1.Range("L2").Value = "'true"
$ws1.Range("L3").Value = "'true"
$ws1.Range("L4").Value = "'false"
$ws1.Range("F5").Value = "'10"
$ws1.Range("L5").Value = "'false"
$ws1.Range("L6").Value = "'false"
$ws1.Range("F7").Value = "'13"
$ws1.Range("L7").Value = "'true"

# --- Browsers (sheet2): new "supported" column (J) ---
$ws2.Range("J1").Value = "supported"
$ws2.Range("J2").Value = "'true"
$ws2.Range("J3").Value = "'true"
$ws2.Range("J4").Value = "'false"
$ws2.Range("J5").Value = "'false"

# --- Selection / active sheet changes ---
$ws2.Range("J10").Select()
$ws1.Activate()
$ws1.Range("K7").Select()
